$wb = $excel.ActiveWorkbook

# --- "combined" sheet (sheet index 1) ---
$ws = $wb.Worksheets.Item(1)
$ws.Range("C13").Value = 6.127418405781572
$ws.Range("N13").Value = 6.124938087420014
$ws.Range("B14").Value = 6.104260631162553
$ws.Range("C14").Value = 6.462910517568996
$ws.Range("D14").Value = 6.104540678677291
$ws.Range("K14").Value = 6.226709285150091
$ws.Range("L14").Value = 6.280524233815036
$ws.Range("M14").Value = 6.319915933542103
$ws.Range("N14").Value = 6.572196768447216
$ws.Range("B15").Value = 6.301891048870691
$ws.Range("C15").Value = 6.785672035018277
$ws.Range("D15").Value = 6.376395790513532
$ws.Range("E15").Value = 6.121230078868197
$ws.Range("J15").Value = 6.200259585932328
$ws.Range("K15").Value = 6.548625524304712
$ws.Range("L15").Value = 6.77349986950787
$ws.Range("M15").Value = 6.756039347987783
$ws.Range("N15").Value = 6.869226906679635
$ws.Range("B16").Value = 6.742501750480489
$ws.Range("D16").Value = 6.738106451659052
$ws.Range("E16").Value = 6.292997314610197
$ws.Range("F16").Value = 6.214665192127313
$ws.Range("G16").Value = 6.072786331362679
$ws.Range("K16").Value = 6.83857291476458
$ws.Range("L16").Value = 7.004091608275051
$ws.Range("M16").Value = 7.10093700329976
$ws.Range("B17").Value = 6.855580305962109
$ws.Range("D17").Value = 7.000313785916195
$ws.Range("F17").Value = 6.420550145818935
$ws.Range("G17").Value = 6.292809252628169
$ws.Range("I17").Value = 6.120603117776967
$ws.Range("J17").Value = 6.7736909842173
$ws.Range("K17").Value = 7.134101028147636
$ws.Range("E18").Value = 6.870591975660952
$ws.Range("H18").Value = 6.155186750824129
$ws.Range("I18").Value = 6.283430704093237
$ws.Range("J18").Value = 7.001403657576703
$ws.Range("F19").Value = 6.978325217049965
$ws.Range("G19").Value = 6.783621649397264
$ws.Range("H19").Value = 6.386330072287121
$ws.Range("G20").Value = 7.021680511385528
$ws.Range("I20").Value = 6.740591213760473
$ws.Range("I21").Value = 6.986575359822184
$ws.Range("H22").Value = 6.812761986345111
$ws.Range("H23").Value = 6.952825668933012
$ws.Range("H24").Value = 7.093340429866335

# --- worksheet index 2 ---
$ws = $wb.Worksheets.Item(2)
$ws.Range("D14").Value = 6.104260631162553
$ws.Range("D15").Value = 6.301891048870691
$ws.Range("D16").Value = 6.742501750480489
$ws.Range("D17").Value = 6.855580305962109

# --- worksheet index 3 ---
$ws = $wb.Worksheets.Item(3)
$ws.Range("D13").Value = 6.127418405781572
$ws.Range("D14").Value = 6.462910517568996
$ws.Range("D15").Value = 6.785672035018277

# --- worksheet index 4 ---
$ws = $wb.Worksheets.Item(4)
$ws.Range("D14").Value = 6.104540678677291
$ws.Range("D15").Value = 6.376395790513532
$ws.Range("D16").Value = 6.738106451659052
$ws.Range("D17").Value = 7.000313785916195

# --- worksheet index 5 ---
$ws = $wb.Worksheets.Item(5)
$ws.Range("D15").Value = 6.121230078868197
$ws.Range("D16").Value = 6.292997314610197
$ws.Range("D18").Value = 6.870591975660952

# --- worksheet index 6 ---
$ws = $wb.Worksheets.Item(6)
$ws.Range("D16").Value = 6.214665192127313
$ws.Range("D17").Value = 6.420550145818935
$ws.Range("D19").Value = 6.978325217049965

# --- worksheet index 7 ---
$ws = $wb.Worksheets.Item(7)
$ws.Range("D16").Value = 6.072786331362679
$ws.Range("D17").Value = 6.292809252628169
$ws.Range("D19").Value = 6.783621649397264
$ws.Range("D20").Value = 7.021680511385528

# --- worksheet index 8 ---
$ws = $wb.Worksheets.Item(8)
$ws.Range("D18").Value = 6.155186750824129
$ws.Range("D19").Value = 6.386330072287121
$ws.Range("D22").Value = 6.812761986345111
$ws.Range("D23").Value = 6.952825668933012
$ws.Range("D24").Value = 7.093340429866335

# --- worksheet index 9 ---
$ws = $wb.Worksheets.Item(9)
$ws.Range("D17").Value = 6.120603117776967
$ws.Range("D18").Value = 6.283430704093237
$ws.Range("D20").Value = 6.740591213760473
$ws.Range("D21").Value = 6.986575359822184

# --- worksheet index 10 ---
$ws = $wb.Worksheets.Item(10)
$ws.Range("D15").Value = 6.200259585932328
$ws.Range("D17").Value = 6.7736909842173
$ws.Range("D18").Value = 7.001403657576703

# --- worksheet index 11 ---
$ws = $wb.Worksheets.Item(11)
$ws.Range("D14").Value = 6.226709285150091
$ws.Range("D15").Value = 6.548625524304712
$ws.Range("D16").Value = 6.83857291476458
$ws.Range("D17").Value = 7.134101028147636

# --- worksheet index 12 ---
$ws = $wb.Worksheets.Item(12)
$ws.Range("D14").Value = 6.280524233815036
$ws.Range("D15").Value = 6.77349986950787
$ws.Range("D16").Value = 7.004091608275051

# --- worksheet index 13 ---
$ws = $wb.Worksheets.Item(13)
$ws.Range("D14").Value = 6.319915933542103
$ws.Range("D15").Value = 6.756039347987783
$ws.Range("D16").Value = 7.10093700329976

# --- worksheet index 14 ---
$ws = $wb.Worksheets.Item(14)
$ws.Range("D13").Value = 6.124938087420014
$ws.Range("D14").Value = 6.572196768447216
$ws.Range("D15").Value = 6.869226906679635
